$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (keep existing style s="1" on A1, extend same style across B1:E1)
$ws.Range("A1").Value = "Volume"
$ws.Range("B1").Value = "Weight [4]"
$ws.Range("C1").Value = "----"
$ws.Range("D1").Value = "----"
$ws.Range("E1").Value = "len_Wi_Hei_Wei_Pack"

$hdr = $ws.Range("B1:E1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Row 2
$ws.Range("A2").Value = 0.24
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = "[]"
$ws.Range("D2").Value = "[]"
$ws.Range("E2").Value = "2x  40x50x60 cm  12 kg/ctn"

# Row 3
$ws.Range("A3").Value = 0.54
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = "[]"
$ws.Range("D3").Value = "[]"
$ws.Range("E3").Value = "12x  50x15x60 cm  17 kg/ctn"

# Row 4
$ws.Range("A4").Value = 6.2208
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = "[]"
$ws.Range("D4").Value = "[]"
$ws.Range("E4").Value = "54x  120x80x12 cm  50 kg/plt"

# Row 5 (totals row: volume/weight only, remaining columns present but empty)
$ws.Range("A5").Value = 7.0008
$ws.Range("B5").Value = 79
$ws.Range("C5").Borders.LineStyle = -4142
$ws.Range("D5").Borders.LineStyle = -4142
$ws.Range("E5").Borders.LineStyle = -4142
